# Updated cryptos list on Tue Jun 20 19:48:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed Price (D) / Volume(1h) (E) values for each coin row.
$updates = @(
    @{ Cell = "D2"; Value = "27.931.18" },
    @{ Cell = "E2"; Value = "  +5.02%  " },
    @{ Cell = "D3"; Value = "1.780.07" },
    @{ Cell = "E3"; Value = "  +3.55%  " },
    @{ Cell = "D4"; Value = "0.9998" },
    @{ Cell = "E4"; Value = "  +0.20%  " },
    @{ Cell = "D5"; Value = "243.39" },
    @{ Cell = "E5"; Value = "  +1.11%  " },
    @{ Cell = "D6"; Value = "0.9997" },
    @{ Cell = "E6"; Value = "  +0.13%  " },
    @{ Cell = "D7"; Value = "0.4890" },
    @{ Cell = "E7"; Value = "  -0.82%  " },
    @{ Cell = "D8"; Value = "0.2663" },
    @{ Cell = "E8"; Value = "  +2.44%  " },
    @{ Cell = "D9"; Value = "0.06246" },
    @{ Cell = "E9"; Value = "  +0.61%  " },
    @{ Cell = "D10"; Value = "1.779.58" },
    @{ Cell = "E10"; Value = "  +3.06%  " },
    @{ Cell = "D11"; Value = "16.30" },
    @{ Cell = "E11"; Value = "  +3.49%  " },
    @{ Cell = "D12"; Value = "0.07006" },
    @{ Cell = "E12"; Value = "  +0.04%  " },
    @{ Cell = "D13"; Value = "0.6223" },
    @{ Cell = "E13"; Value = "  +2.59%  " },
    @{ Cell = "D14"; Value = "4.618" },
    @{ Cell = "E14"; Value = "  +2.89%  " },
    @{ Cell = "D15"; Value = "79.77" },
    @{ Cell = "E15"; Value = "  +3.91%  " },
    @{ Cell = "D16"; Value = "1.000" },
    @{ Cell = "D17"; Value = "27.900.33" },
    @{ Cell = "E17"; Value = "  +5.51%  " },
    @{ Cell = "D18"; Value = "0.9994" },
    @{ Cell = "E18"; Value = "  +0.18%  " },
    @{ Cell = "D19"; Value = "0.000007200" },
    @{ Cell = "E19"; Value = "  +0.70%  " },
    @{ Cell = "D20"; Value = "11.83" },
    @{ Cell = "E20"; Value = "  +4.13%  " },
    @{ Cell = "D21"; Value = "2.008.02" },
    @{ Cell = "E21"; Value = "  +3.35%  " },
    @{ Cell = "D22"; Value = "4.570" },
    @{ Cell = "D23"; Value = "8.660" },
    @{ Cell = "E23"; Value = "  +1.71%  " },
    @{ Cell = "D24"; Value = "5.208" },
    @{ Cell = "E24"; Value = "  +2.48%  " },
    @{ Cell = "D25"; Value = "141.78" },
    @{ Cell = "E25"; Value = "  +2.90%  " },
    @{ Cell = "D26"; Value = "15.57" },
    @{ Cell = "E26"; Value = "  +2.04%  " },
    @{ Cell = "D27"; Value = "1.861" },
    @{ Cell = "E27"; Value = "  +7.17%  " },
    @{ Cell = "D28"; Value = "108.86" },
    @{ Cell = "E28"; Value = "  +2.98%  " },
    @{ Cell = "D29"; Value = "1.386" },
    @{ Cell = "E29"; Value = "  -1.17%  " },
    @{ Cell = "D30"; Value = "4.158" },
    @{ Cell = "E30"; Value = "  +6.09%  " },
    @{ Cell = "D31"; Value = "0.08233" },
    @{ Cell = "E31"; Value = "  +3.51%  " },
    @{ Cell = "D32"; Value = "3.789" },
    @{ Cell = "E32"; Value = "  +3.90%  " },
    @{ Cell = "D33"; Value = "0.04759" },
    @{ Cell = "E33"; Value = "  +5.86%  " },
    @{ Cell = "D34"; Value = "1.068" },
    @{ Cell = "E34"; Value = "  +6.90%  " },
    @{ Cell = "D35"; Value = "2.601" },
    @{ Cell = "E35"; Value = "  -0.36%  " },
    @{ Cell = "D36"; Value = "0.6439" },
    @{ Cell = "E36"; Value = "  +3.32%  " },
    @{ Cell = "D37"; Value = "0.9440" },
    @{ Cell = "E37"; Value = "  +0.78%  " },
    @{ Cell = "D38"; Value = "2.585" },
    @{ Cell = "E38"; Value = "  +7.25%  " },
    @{ Cell = "D39"; Value = "2.047" },
    @{ Cell = "E39"; Value = "  +2.44%  " },
    @{ Cell = "D40"; Value = "5.912" },
    @{ Cell = "E40"; Value = "  +7.21%  " },
    @{ Cell = "D41"; Value = "0.01539" },
    @{ Cell = "E41"; Value = "  +1.82%  " },
    @{ Cell = "D42"; Value = "0.9996" },
    @{ Cell = "E42"; Value = "  +0.13%  " },
    @{ Cell = "D43"; Value = "99.98" },
    @{ Cell = "E43"; Value = "  +0.80%  " },
    @{ Cell = "D44"; Value = "0.3958" },
    @{ Cell = "E44"; Value = "  +3.19%  " },
    @{ Cell = "D45"; Value = "7.217" },
    @{ Cell = "E45"; Value = "  +4.45%  " },
    @{ Cell = "E46"; Value = "  +3.56%  " },
    @{ Cell = "E47"; Value = "  +0.69%  " },
    @{ Cell = "D48"; Value = "7.958" },
    @{ Cell = "E48"; Value = "  +3.04%  " },
    @{ Cell = "D49"; Value = "1.285" },
    @{ Cell = "E49"; Value = "  +4.91%  " },
    @{ Cell = "D50"; Value = "30.46" },
    @{ Cell = "E50"; Value = "  +1.19%  " },
    @{ Cell = "D51"; Value = "52.60" },
    @{ Cell = "E51"; Value = "  +2.34%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Leading apostrophe forces Excel to store the new value as literal
    # text (matching the original cell content type) instead of letting
    # automatic type inference coerce number-like strings (e.g. "1.000",
    # "16.30", "0.000007200") into numeric values and losing formatting.
    $cell.Value = "'" + $u.Value
    # Restore the default "Normal" style so forcing text entry above does
    # not leave a stray number-format/quote-prefix style on the cell.
    $cell.Style = "Normal"
}
